$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 02:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 923812
$ws.Range("C4").Value = 37370
$ws.Range("D4").Value = 110400
$ws.Range("E4").Value = 761315
$ws.Range("F4").Value = 15097
$ws.Range("G4").Value = 1863
$ws.Range("H4").Value = 52097

# --- Row 30: Singapur ---
$ws.Range("D30").Value = 956
$ws.Range("E30").Value = 11107

# --- Row 57: Argentina ---
$ws.Range("B57").Value = 3607
$ws.Range("C57").Value = 172
$ws.Range("E57").Value = 2455
$ws.Range("G57").Value = 11
$ws.Range("H57").Value = 176

# --- Row 93: Principado de Andorra ---
$ws.Range("D93").Value = 344
$ws.Range("E93").Value = 347

# --- Row 154: Barbados ---
$ws.Range("B154").Value = 77
$ws.Range("C154").Value = 1
$ws.Range("D154").Value = 31

# --- Rows 179/180: Antigua y Barbuda and Timor Oriental swap order,
# with Timor Oriental's data refreshed ---
$ws.Range("A179").Value = "Timor Oriental"
$ws.Range("B179").Value = 24
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 2
$ws.Range("E179").Value = 22
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

$ws.Range("A180").Value = "Antigua y Barbuda"
$ws.Range("B180").Value = 24
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 10
$ws.Range("E180").Value = 11
$ws.Range("F180").Value = 1
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 3

# --- Row 202: Surinam ---
$ws.Range("D202").Value = 7
$ws.Range("E202").Value = 2
